# Update the header date and the 100 arithmetic answer-sheet cells to the
# new set of equations (see commit "Update master to output generated at
# 2dcbd77"). Each old string is unique in the document, so a plain
# Find/Replace (MatchCase, not whole-word since strings contain +/-/=)
# against the whole document content is sufficient and exact.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-03-21 Tuesday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-03-22 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("3+82=85", $true, $true, $false, $false, $false, $true, 1, $false, "54-21=33", 2) | Out-Null
$d.Content.Find.Execute("83-74=9", $true, $true, $false, $false, $false, $true, 1, $false, "92-54=38", 2) | Out-Null
$d.Content.Find.Execute("50-12=38", $true, $true, $false, $false, $false, $true, 1, $false, "75-48=27", 2) | Out-Null
$d.Content.Find.Execute("83-50=33", $true, $true, $false, $false, $false, $true, 1, $false, "84-82=2", 2) | Out-Null
$d.Content.Find.Execute("27+70=97", $true, $true, $false, $false, $false, $true, 1, $false, "57-24=33", 2) | Out-Null
$d.Content.Find.Execute("38+61=99", $true, $true, $false, $false, $false, $true, 1, $false, "53+41=94", 2) | Out-Null
$d.Content.Find.Execute("63+5=68", $true, $true, $false, $false, $false, $true, 1, $false, "39-25=14", 2) | Out-Null
$d.Content.Find.Execute("72+6=78", $true, $true, $false, $false, $false, $true, 1, $false, "37+62=99", 2) | Out-Null
$d.Content.Find.Execute("80-33=47", $true, $true, $false, $false, $false, $true, 1, $false, "73-38=35", 2) | Out-Null
$d.Content.Find.Execute("93-11=82", $true, $true, $false, $false, $false, $true, 1, $false, "4+38=42", 2) | Out-Null
$d.Content.Find.Execute("17+15=32", $true, $true, $false, $false, $false, $true, 1, $false, "51-19=32", 2) | Out-Null
$d.Content.Find.Execute("44-40=4", $true, $true, $false, $false, $false, $true, 1, $false, "99-40=59", 2) | Out-Null
$d.Content.Find.Execute("30+32=62", $true, $true, $false, $false, $false, $true, 1, $false, "86-24=62", 2) | Out-Null
$d.Content.Find.Execute("23+20=43", $true, $true, $false, $false, $false, $true, 1, $false, "62-30=32", 2) | Out-Null
$d.Content.Find.Execute("84-21=63", $true, $true, $false, $false, $false, $true, 1, $false, "69+29=98", 2) | Out-Null
$d.Content.Find.Execute("36+45=81", $true, $true, $false, $false, $false, $true, 1, $false, "75-71=4", 2) | Out-Null
$d.Content.Find.Execute("9+56=65", $true, $true, $false, $false, $false, $true, 1, $false, "16+66=82", 2) | Out-Null
$d.Content.Find.Execute("0+89=89", $true, $true, $false, $false, $false, $true, 1, $false, "70+7=77", 2) | Out-Null
$d.Content.Find.Execute("57+11=68", $true, $true, $false, $false, $false, $true, 1, $false, "41+15=56", 2) | Out-Null
$d.Content.Find.Execute("94-9=85", $true, $true, $false, $false, $false, $true, 1, $false, "95-77=18", 2) | Out-Null
$d.Content.Find.Execute("48-14=34", $true, $true, $false, $false, $false, $true, 1, $false, "1+2=3", 2) | Out-Null
$d.Content.Find.Execute("58-46=12", $true, $true, $false, $false, $false, $true, 1, $false, "97-74=23", 2) | Out-Null
$d.Content.Find.Execute("2+24=26", $true, $true, $false, $false, $false, $true, 1, $false, "43-32=11", 2) | Out-Null
$d.Content.Find.Execute("17+6=23", $true, $true, $false, $false, $false, $true, 1, $false, "49+12=61", 2) | Out-Null
$d.Content.Find.Execute("43+41=84", $true, $true, $false, $false, $false, $true, 1, $false, "13+8=21", 2) | Out-Null
$d.Content.Find.Execute("48+38=86", $true, $true, $false, $false, $false, $true, 1, $false, "0+1=1", 2) | Out-Null
$d.Content.Find.Execute("25+36=61", $true, $true, $false, $false, $false, $true, 1, $false, "19-17=2", 2) | Out-Null
$d.Content.Find.Execute("96-32=64", $true, $true, $false, $false, $false, $true, 1, $false, "97-70=27", 2) | Out-Null
$d.Content.Find.Execute("87-42=45", $true, $true, $false, $false, $false, $true, 1, $false, "91-81=10", 2) | Out-Null
$d.Content.Find.Execute("26+24=50", $true, $true, $false, $false, $false, $true, 1, $false, "94-43=51", 2) | Out-Null
$d.Content.Find.Execute("34-26=8", $true, $true, $false, $false, $false, $true, 1, $false, "95-66=29", 2) | Out-Null
$d.Content.Find.Execute("49-35=14", $true, $true, $false, $false, $false, $true, 1, $false, "98-40=58", 2) | Out-Null
$d.Content.Find.Execute("90-12=78", $true, $true, $false, $false, $false, $true, 1, $false, "66+4=70", 2) | Out-Null
$d.Content.Find.Execute("90-84=6", $true, $true, $false, $false, $false, $true, 1, $false, "29+21=50", 2) | Out-Null
$d.Content.Find.Execute("84+12=96", $true, $true, $false, $false, $false, $true, 1, $false, "7+34=41", 2) | Out-Null
$d.Content.Find.Execute("66-64=2", $true, $true, $false, $false, $false, $true, 1, $false, "13+86=99", 2) | Out-Null
$d.Content.Find.Execute("85-72=13", $true, $true, $false, $false, $false, $true, 1, $false, "68-27=41", 2) | Out-Null
$d.Content.Find.Execute("3+54=57", $true, $true, $false, $false, $false, $true, 1, $false, "21-5=16", 2) | Out-Null
$d.Content.Find.Execute("51-12=39", $true, $true, $false, $false, $false, $true, 1, $false, "47-18=29", 2) | Out-Null
$d.Content.Find.Execute("42+49=91", $true, $true, $false, $false, $false, $true, 1, $false, "53-15=38", 2) | Out-Null
$d.Content.Find.Execute("16+40=56", $true, $true, $false, $false, $false, $true, 1, $false, "74-10=64", 2) | Out-Null
$d.Content.Find.Execute("2+38=40", $true, $true, $false, $false, $false, $true, 1, $false, "80+13=93", 2) | Out-Null
$d.Content.Find.Execute("32+40=72", $true, $true, $false, $false, $false, $true, 1, $false, "52+8=60", 2) | Out-Null
$d.Content.Find.Execute("50+4=54", $true, $true, $false, $false, $false, $true, 1, $false, "5+88=93", 2) | Out-Null
$d.Content.Find.Execute("79+17=96", $true, $true, $false, $false, $false, $true, 1, $false, "7-2=5", 2) | Out-Null
$d.Content.Find.Execute("93-84=9", $true, $true, $false, $false, $false, $true, 1, $false, "53+27=80", 2) | Out-Null
$d.Content.Find.Execute("54-48=6", $true, $true, $false, $false, $false, $true, 1, $false, "36-30=6", 2) | Out-Null
$d.Content.Find.Execute("11+80=91", $true, $true, $false, $false, $false, $true, 1, $false, "16+9=25", 2) | Out-Null
$d.Content.Find.Execute("39-9=30", $true, $true, $false, $false, $false, $true, 1, $false, "63+16=79", 2) | Out-Null
$d.Content.Find.Execute("39+10=49", $true, $true, $false, $false, $false, $true, 1, $false, "15+0=15", 2) | Out-Null
$d.Content.Find.Execute("3+75=78", $true, $true, $false, $false, $false, $true, 1, $false, "97-40=57", 2) | Out-Null
$d.Content.Find.Execute("52+15=67", $true, $true, $false, $false, $false, $true, 1, $false, "27-22=5", 2) | Out-Null
$d.Content.Find.Execute("14+13=27", $true, $true, $false, $false, $false, $true, 1, $false, "71-59=12", 2) | Out-Null
$d.Content.Find.Execute("12+70=82", $true, $true, $false, $false, $false, $true, 1, $false, "90-21=69", 2) | Out-Null
$d.Content.Find.Execute("69-59=10", $true, $true, $false, $false, $false, $true, 1, $false, "70-7=63", 2) | Out-Null
$d.Content.Find.Execute("74-4=70", $true, $true, $false, $false, $false, $true, 1, $false, "48+34=82", 2) | Out-Null
$d.Content.Find.Execute("71-43=28", $true, $true, $false, $false, $false, $true, 1, $false, "76-41=35", 2) | Out-Null
$d.Content.Find.Execute("69-39=30", $true, $true, $false, $false, $false, $true, 1, $false, "34+54=88", 2) | Out-Null
$d.Content.Find.Execute("46-7=39", $true, $true, $false, $false, $false, $true, 1, $false, "84-15=69", 2) | Out-Null
$d.Content.Find.Execute("3+33=36", $true, $true, $false, $false, $false, $true, 1, $false, "2+41=43", 2) | Out-Null
$d.Content.Find.Execute("68-11=57", $true, $true, $false, $false, $false, $true, 1, $false, "69-36=33", 2) | Out-Null
$d.Content.Find.Execute("59-32=27", $true, $true, $false, $false, $false, $true, 1, $false, "81-12=69", 2) | Out-Null
$d.Content.Find.Execute("78-46=32", $true, $true, $false, $false, $false, $true, 1, $false, "22+76=98", 2) | Out-Null
$d.Content.Find.Execute("4+66=70", $true, $true, $false, $false, $false, $true, 1, $false, "86-4=82", 2) | Out-Null
$d.Content.Find.Execute("17+62=79", $true, $true, $false, $false, $false, $true, 1, $false, "73-6=67", 2) | Out-Null
$d.Content.Find.Execute("26+25=51", $true, $true, $false, $false, $false, $true, 1, $false, "35-7=28", 2) | Out-Null
$d.Content.Find.Execute("50-17=33", $true, $true, $false, $false, $false, $true, 1, $false, "77-70=7", 2) | Out-Null
$d.Content.Find.Execute("61+30=91", $true, $true, $false, $false, $false, $true, 1, $false, "20+44=64", 2) | Out-Null
$d.Content.Find.Execute("68+5=73", $true, $true, $false, $false, $false, $true, 1, $false, "90-57=33", 2) | Out-Null
$d.Content.Find.Execute("83-5=78", $true, $true, $false, $false, $false, $true, 1, $false, "92-56=36", 2) | Out-Null
$d.Content.Find.Execute("87-62=25", $true, $true, $false, $false, $false, $true, 1, $false, "25+67=92", 2) | Out-Null
$d.Content.Find.Execute("80-50=30", $true, $true, $false, $false, $false, $true, 1, $false, "96-6=90", 2) | Out-Null
$d.Content.Find.Execute("30+25=55", $true, $true, $false, $false, $false, $true, 1, $false, "59+14=73", 2) | Out-Null
$d.Content.Find.Execute("86-26=60", $true, $true, $false, $false, $false, $true, 1, $false, "76-75=1", 2) | Out-Null
$d.Content.Find.Execute("88+0=88", $true, $true, $false, $false, $false, $true, 1, $false, "94-29=65", 2) | Out-Null
$d.Content.Find.Execute("68-40=28", $true, $true, $false, $false, $false, $true, 1, $false, "90-18=72", 2) | Out-Null
$d.Content.Find.Execute("99-88=11", $true, $true, $false, $false, $false, $true, 1, $false, "79-42=37", 2) | Out-Null
$d.Content.Find.Execute("25+29=54", $true, $true, $false, $false, $false, $true, 1, $false, "60-45=15", 2) | Out-Null
$d.Content.Find.Execute("86+7=93", $true, $true, $false, $false, $false, $true, 1, $false, "15+31=46", 2) | Out-Null
$d.Content.Find.Execute("6+65=71", $true, $true, $false, $false, $false, $true, 1, $false, "58+41=99", 2) | Out-Null
$d.Content.Find.Execute("63+29=92", $true, $true, $false, $false, $false, $true, 1, $false, "39+25=64", 2) | Out-Null
$d.Content.Find.Execute("93-65=28", $true, $true, $false, $false, $false, $true, 1, $false, "82-30=52", 2) | Out-Null
$d.Content.Find.Execute("5+70=75", $true, $true, $false, $false, $false, $true, 1, $false, "11+21=32", 2) | Out-Null
$d.Content.Find.Execute("31-30=1", $true, $true, $false, $false, $false, $true, 1, $false, "65-21=44", 2) | Out-Null
$d.Content.Find.Execute("91+0=91", $true, $true, $false, $false, $false, $true, 1, $false, "12+56=68", 2) | Out-Null
$d.Content.Find.Execute("79-57=22", $true, $true, $false, $false, $false, $true, 1, $false, "30+7=37", 2) | Out-Null
$d.Content.Find.Execute("21+50=71", $true, $true, $false, $false, $false, $true, 1, $false, "45+8=53", 2) | Out-Null
$d.Content.Find.Execute("51-23=28", $true, $true, $false, $false, $false, $true, 1, $false, "73-50=23", 2) | Out-Null
$d.Content.Find.Execute("99-77=22", $true, $true, $false, $false, $false, $true, 1, $false, "46-11=35", 2) | Out-Null
$d.Content.Find.Execute("14+57=71", $true, $true, $false, $false, $false, $true, 1, $false, "99-60=39", 2) | Out-Null
$d.Content.Find.Execute("18-18=0", $true, $true, $false, $false, $false, $true, 1, $false, "21+41=62", 2) | Out-Null
$d.Content.Find.Execute("33-6=27", $true, $true, $false, $false, $false, $true, 1, $false, "18+39=57", 2) | Out-Null
$d.Content.Find.Execute("36-25=11", $true, $true, $false, $false, $false, $true, 1, $false, "71-24=47", 2) | Out-Null
$d.Content.Find.Execute("25+60=85", $true, $true, $false, $false, $false, $true, 1, $false, "33+6=39", 2) | Out-Null
$d.Content.Find.Execute("33+4=37", $true, $true, $false, $false, $false, $true, 1, $false, "62-27=35", 2) | Out-Null
$d.Content.Find.Execute("58+13=71", $true, $true, $false, $false, $false, $true, 1, $false, "2+64=66", 2) | Out-Null
$d.Content.Find.Execute("10+55=65", $true, $true, $false, $false, $false, $true, 1, $false, "42-0=42", 2) | Out-Null
$d.Content.Find.Execute("79-20=59", $true, $true, $false, $false, $false, $true, 1, $false, "58-47=11", 2) | Out-Null
$d.Content.Find.Execute("18+16=34", $true, $true, $false, $false, $false, $true, 1, $false, "44-30=14", 2) | Out-Null
$d.Content.Find.Execute("6+17=23", $true, $true, $false, $false, $false, $true, 1, $false, "75-46=29", 2) | Out-Null
